$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its text (inline-string-like) representation even
# for values that look numeric (e.g. "1.00", "0.0280", "239.91") so trailing
# zeros / grouping dots are preserved exactly as in the source data.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "43.807.53"
$ws.Range("E2").Value = "  -0.91%  "

# Row 3
$ws.Range("D3").Value = "2.349.09"
$ws.Range("E3").Value = "  -0.25%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("E5").Value = "  -0.65%  "

# Row 6
$ws.Range("D6").Value = "239.91"
$ws.Range("E6").Value = "  -0.26%  "

# Row 7
$ws.Range("D7").Value = "73.24"
$ws.Range("E7").Value = "  -0.64%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("D9").Value = "0.601"
$ws.Range("E9").Value = "  +8.03%  "

# Row 10
$ws.Range("E10").Value = "  -2.51%  "

# Row 11
$ws.Range("D11").Value = "58.72"
$ws.Range("E11").Value = "  +2.79%  "

# Row 12
$ws.Range("D12").Value = "32.79"
$ws.Range("E12").Value = "  +6.08%  "

# Row 13
$ws.Range("D13").Value = "7.31"
$ws.Range("E13").Value = "  -1.07%  "

# Row 14
$ws.Range("E14").Value = "  +0.13%  "

# Row 15
$ws.Range("D15").Value = "2.699.65"
$ws.Range("E15").Value = "  -0.22%  "

# Row 16
$ws.Range("D16").Value = "16.43"
$ws.Range("E16").Value = "  -2.30%  "

# Row 17
$ws.Range("D17").Value = "0.904"
$ws.Range("E17").Value = "  -1.01%  "

# Row 18
$ws.Range("D18").Value = "2.354.34"
$ws.Range("E18").Value = "  -0.01%  "

# Row 19
$ws.Range("D19").Value = "43.714.52"
$ws.Range("E19").Value = "  -1.22%  "

# Row 20
$ws.Range("E20").Value = "  -1.02%  "

# Row 21
$ws.Range("D21").Value = "6.75"
$ws.Range("E21").Value = "  +1.51%  "

# Row 22
$ws.Range("D22").Value = "77.43"
$ws.Range("E22").Value = "  -0.50%  "

# Row 23
$ws.Range("D23").Value = "256.87"
$ws.Range("E23").Value = "  +0.28%  "

# Row 24
$ws.Range("D24").Value = "2.03"
$ws.Range("E24").Value = "  +28.32%  "

# Row 25
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.07%  "

# Row 26
$ws.Range("D26").Value = "3.73"
$ws.Range("E26").Value = "  -1.35%  "

# Row 27
$ws.Range("E27").Value = "  -2.21%  "

# Row 28
$ws.Range("D28").Value = "10.62"
$ws.Range("E28").Value = "  -0.25%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "22.68"
$ws.Range("E29").Value = "  +0.16%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.21"
$ws.Range("E30").Value = "  -4.08%  "

# Row 31
$ws.Range("D31").Value = "175.80"
$ws.Range("E31").Value = "  +0.69%  "

# Row 32
$ws.Range("D32").Value = "0.131"
$ws.Range("E32").Value = "  -0.91%  "

# Row 33
$ws.Range("D33").Value = "0.137"
$ws.Range("E33").Value = "  +2.91%  "

# Row 34
$ws.Range("E34").Value = "  +1.77%  "

# Row 35
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "5.18"
$ws.Range("E35").Value = "  -3.87%  "

# Row 36
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "5.47"
$ws.Range("E36").Value = "  +2.50%  "

# Row 37
$ws.Range("D37").Value = "3.80"
$ws.Range("E37").Value = "  -3.07%  "

# Row 38
$ws.Range("D38").Value = "2.36"

# Row 39
$ws.Range("B39").Value = "THORChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D39").Value = "6.26"
$ws.Range("E39").Value = "  -4.62%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.0280"
$ws.Range("E40").Value = "  +2.41%  "

# Row 41
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").Value = "69.33"
$ws.Range("E41").Value = "  +31.22%  "

# Row 42
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").Value = "0.111"
$ws.Range("E42").Value = "  +10.85%  "

# Row 43
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "0.204"
$ws.Range("E43").Value = "  +9.49%  "

# Row 44
$ws.Range("D44").Value = "9.08"
$ws.Range("E44").Value = "  +0.25%  "

# Row 45
$ws.Range("D45").Value = "18.93"
$ws.Range("E45").Value = "  -1.59%  "

# Row 46
$ws.Range("D46").Value = "4.78"
$ws.Range("E46").Value = "  +6.48%  "

# Row 47
$ws.Range("D47").Value = "2.53"
$ws.Range("E47").Value = "  +3.22%  "

# Row 48
$ws.Range("E48").Value = "  +0.08%  "

# Row 49
$ws.Range("D49").Value = "1.24"
$ws.Range("E49").Value = "  -1.52%  "

# Row 50
$ws.Range("D50").Value = "99.47"
$ws.Range("E50").Value = "  -1.00%  "

# Row 51
$ws.Range("D51").Value = "1.16"
$ws.Range("E51").Value = "  -0.97%  "
